# Stoppage Report User Interface
# Update the single data row (row 2) on the active worksheet to reflect
# the newly generated stoppage report: only "Gibson Ring Stops" (column E)
# has a stoppage (1, 100%), all other stop counters are 0 (0%), and both
# the report start/end dates are the same day, 2017-04-07.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "0 (0%)"
$ws.Range("B2").Value = "0 (0%)"
$ws.Range("C2").Value = "0 (0%)"
$ws.Range("D2").Value = "0 (0%)"
$ws.Range("E2").Value = "1 (100%)"
$ws.Range("F2").Value = "0 (0%)"
$ws.Range("G2").Value = "0 (0%)"
$ws.Range("H2").Value = "0 (0%)"
$ws.Range("I2").Value = "0 (0%)"
$ws.Range("J2").Value = "0 (0%)"
$ws.Range("K2").Value = "'2017-04-07"
$ws.Range("L2").Value = "'2017-04-07"
